# Update test data for Recommended Content:
# The cardImageSrc path no longer contains the dated "2019-11/" subfolder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_recommended_content")

$newPath = "/sites/default/files/styles/cgov_featured/public/cgov_image/featured/"

$ws.Range("J2").Value = $newPath
$ws.Range("J3").Value = $newPath
$ws.Range("J4").Value = $newPath

# Leave the selection where the author left it when saving.
$ws.Range("J7").Select() | Out-Null
